$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(64, 8).Value = 3866.6667
$ws.Cells.Item(64, 9).Value = 3800
$ws.Cells.Item(64, 10).Value = 4000
$ws.Cells.Item(64, 11).Value = 3800
$ws.Cells.Item(64, 12).Value = 4000
$ws.Cells.Item(64, 13).Value = -3552
$ws.Cells.Item(64, 14).Value = -4496

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(67, 8).Value = 3866.6667
$ws.Cells.Item(67, 9).Value = 3800
$ws.Cells.Item(67, 10).Value = 4000
$ws.Cells.Item(67, 11).Value = 3800
$ws.Cells.Item(67, 12).Value = 4000
$ws.Cells.Item(67, 13).Value = -2942
$ws.Cells.Item(67, 14).Value = -5716

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(113, 8).Value = 2152.5
$ws.Cells.Item(113, 9).Value = 2251.25
$ws.Cells.Item(113, 10).Value = 2130.5557
$ws.Cells.Item(113, 11).Value = 2251.25
$ws.Cells.Item(113, 12).Value = 2130.5557
$ws.Cells.Item(113, 13).Value = 1002.75
$ws.Cells.Item(113, 14).Value = -8638.555700000001

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(116, 8).Value = 2708.5715
$ws.Cells.Item(116, 9).Value = 2515
$ws.Cells.Item(116, 10).Value = 2966.6667
$ws.Cells.Item(116, 11).Value = 2515
$ws.Cells.Item(116, 12).Value = 2966.6667
$ws.Cells.Item(116, 13).Value = 927
$ws.Cells.Item(116, 14).Value = -9850.6667

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(132, 8).Value = 9264802
$ws.Cells.Item(132, 9).Value = 11115878
$ws.Cells.Item(132, 11).Value = 33347634
$ws.Cells.Item(132, 13).Value = -33345104

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 8773106
$ws.Cells.Item(61, 9).Value = 9260444
$ws.Cells.Item(61, 10).Value = 1014
$ws.Cells.Item(61, 11).Value = 9260444
$ws.Cells.Item(61, 12).Value = 1014
$ws.Cells.Item(61, 13).Value = -9260232
$ws.Cells.Item(61, 14).Value = -1438

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(132, 8).Value = 2102533.8
$ws.Cells.Item(132, 9).Value = 1360.2106
$ws.Cells.Item(132, 10).Value = 6538345
$ws.Cells.Item(132, 11).Value = 4080.6318
$ws.Cells.Item(132, 12).Value = 19615035
$ws.Cells.Item(132, 13).Value = -1550.6318
$ws.Cells.Item(132, 14).Value = -19620095

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value = 8773106
$ws.Cells.Item(136, 9).Value = 9260444
$ws.Cells.Item(136, 10).Value = 1014
$ws.Cells.Item(136, 11).Value = 27781332
$ws.Cells.Item(136, 12).Value = 3042
$ws.Cells.Item(136, 13).Value = -27778782
$ws.Cells.Item(136, 14).Value = -8142

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(86, 8).Value = 752366.6
$ws.Cells.Item(86, 9).Value = 2477.1428
$ws.Cells.Item(86, 10).Value = 1369922.8
$ws.Cells.Item(86, 11).Value = 2477.1428
$ws.Cells.Item(86, 12).Value = 1369922.8
$ws.Cells.Item(86, 13).Value = -1354.1428
$ws.Cells.Item(86, 14).Value = -1372168.8

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(89, 8).Value = 752366.6
$ws.Cells.Item(89, 9).Value = 2477.1428
$ws.Cells.Item(89, 10).Value = 1369922.8
$ws.Cells.Item(89, 11).Value = 12385.714
$ws.Cells.Item(89, 12).Value = 6849614
$ws.Cells.Item(89, 13).Value = -6769.714
$ws.Cells.Item(89, 14).Value = -6860846

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(60, 8).Value = 25235.334
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 25235.334
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 12).Value = 25235.334
$ws.Cells.Item(60, 13).ClearContents()
$ws.Cells.Item(60, 14).Value = -26257.334

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(99, 8).Value = 55558860
$ws.Cells.Item(99, 9).Value = 111113600
$ws.Cells.Item(99, 10).Value = 4122.222
$ws.Cells.Item(99, 11).Value = 111113600
$ws.Cells.Item(99, 12).Value = 4122.222
$ws.Cells.Item(99, 13).Value = -111112102
$ws.Cells.Item(99, 14).Value = -7118.222

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(126, 8).Value = 55558860
$ws.Cells.Item(126, 9).Value = 111113600
$ws.Cells.Item(126, 10).Value = 4122.222
$ws.Cells.Item(126, 11).Value = 333340800
$ws.Cells.Item(126, 12).Value = 12366.666
$ws.Cells.Item(126, 13).Value = -333338330
$ws.Cells.Item(126, 14).Value = -17306.666

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(132, 8).Value = 10102184
$ws.Cells.Item(132, 9).Value = 1048.4783
$ws.Cells.Item(132, 10).Value = 33334794
$ws.Cells.Item(132, 11).Value = 3145.4349
$ws.Cells.Item(132, 12).Value = 100004382
$ws.Cells.Item(132, 13).Value = -615.4349000000002
$ws.Cells.Item(132, 14).Value = -100009442

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value = 797.0599999999999
$ws.Cells.Item(134, 9).Value = 1113.65
$ws.Cells.Item(134, 10).Value = 586
$ws.Cells.Item(134, 11).Value = 3340.95
$ws.Cells.Item(134, 12).Value = 1758
$ws.Cells.Item(134, 13).Value = -805.9500000000003
$ws.Cells.Item(134, 14).Value = -6828

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(113, 8).Value = 46970096
$ws.Cells.Item(113, 9).Value = 27778172
$ws.Cells.Item(113, 10).Value = 54167064
$ws.Cells.Item(113, 11).Value = 83334516
$ws.Cells.Item(113, 12).Value = 162501192
$ws.Cells.Item(113, 13).Value = -83332346
$ws.Cells.Item(113, 14).Value = -162505532

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(129, 8).Value = 1505.625
$ws.Cells.Item(129, 9).Value = 1522.5
$ws.Cells.Item(129, 10).Value = 1500
$ws.Cells.Item(129, 11).Value = 4567.5
$ws.Cells.Item(129, 12).Value = 4500
$ws.Cells.Item(129, 13).Value = 432.5
$ws.Cells.Item(129, 14).Value = -14500

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(131, 8).Value = 877.0612
$ws.Cells.Item(131, 10).Value = 898.5161000000001
$ws.Cells.Item(131, 12).Value = 2695.5483
$ws.Cells.Item(131, 14).Value = -12775.5483

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).ClearContents()

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(80, 8).Value = 4548295.5
$ws.Cells.Item(80, 9).Value = 3508.3333
$ws.Cells.Item(80, 10).Value = 10002040
$ws.Cells.Item(80, 11).Value = 3508.3333
$ws.Cells.Item(80, 12).Value = 10002040
$ws.Cells.Item(80, 13).Value = -2510.3333
$ws.Cells.Item(80, 14).Value = -10004036

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(83, 8).Value = 4548295.5
$ws.Cells.Item(83, 9).Value = 3508.3333
$ws.Cells.Item(83, 10).Value = 10002040
$ws.Cells.Item(83, 11).Value = 17541.6665
$ws.Cells.Item(83, 12).Value = 50010200
$ws.Cells.Item(83, 13).Value = -12549.6665
$ws.Cells.Item(83, 14).Value = -50020184

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value = 10912.333
$ws.Cells.Item(132, 9).Value = 2327.7778
$ws.Cells.Item(132, 10).Value = 36666
$ws.Cells.Item(132, 11).Value = 6983.3334
$ws.Cells.Item(132, 12).Value = 109998
$ws.Cells.Item(132, 13).Value = -4453.3334
$ws.Cells.Item(132, 14).Value = -115058

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(140, 8).Value = 48500
$ws.Cells.Item(140, 10).Value = 48500
$ws.Cells.Item(140, 12).Value = 48500
$ws.Cells.Item(140, 14).Value = -58860

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(141, 8).Value = 48107.25
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 48107.25
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 48107.25
$ws.Cells.Item(141, 13).ClearContents()
$ws.Cells.Item(141, 14).Value = -58467.25

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(61, 8).Value = 1372.1578
$ws.Cells.Item(61, 9).Value = 1614.5
$ws.Cells.Item(61, 10).Value = 1260.3077
$ws.Cells.Item(61, 11).Value = 1614.5
$ws.Cells.Item(61, 12).Value = 1260.3077
$ws.Cells.Item(61, 13).Value = -1412.5
$ws.Cells.Item(61, 14).Value = -1664.3077

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(113, 8).Value = 1372.1578
$ws.Cells.Item(113, 9).Value = 1614.5
$ws.Cells.Item(113, 10).Value = 1260.3077
$ws.Cells.Item(113, 11).Value = 1614.5
$ws.Cells.Item(113, 12).Value = 1260.3077
$ws.Cells.Item(113, 13).Value = 555.5
$ws.Cells.Item(113, 14).Value = -5600.3077

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(132, 8).Value = 10589.725
$ws.Cells.Item(132, 9).Value = 2775.2666
$ws.Cells.Item(132, 11).Value = 8325.799800000001
$ws.Cells.Item(132, 13).Value = -5795.799800000001

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(113, 8).Value = 1063
$ws.Cells.Item(113, 9).Value = 1063
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 3189
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -1019
$ws.Cells.Item(113, 14).ClearContents()
